$wb = $excel.ActiveWorkbook

# Update "展览" sheet (sheet1) - "想去人数" (want-to-go count) values
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F6").Value = 465
$ws1.Range("F9").Value = 610

# Update "全部类型" sheet (sheet4) - mirrors the same rows
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F6").Value = 465
$ws4.Range("F9").Value = 610
